# budget CA test cases
# - Anonymous sheet (sheet1): column A used to be tagged with the
#   "AnonymousUserTestCases" label on every data row; that label now only
#   belongs on the new CA test rows living on Sheet1, so clear it here.
# - Sheet1 (sheet2): rows 2-15 are the new CA test cases, tag column A with
#   the "AnonymousUserTestCases" label, add an explicit width for column A,
#   and fill in the Airline/FlightNumber columns for row 12.
# - Active window/tab moves from the Anonymous sheet to Sheet1.

$wb = $excel.ActiveWorkbook
$wsAnonymous = $wb.Worksheets.Item(1)
$wsSheet1 = $wb.Worksheets.Item(2)

# --- Anonymous sheet: drop the per-row "AnonymousUserTestCases" tag in column A ---
$wsAnonymous.Range("A2:A21").Value = ""

# --- Sheet1: tag column A with "AnonymousUserTestCases" for the new CA rows ---
for ($r = 2; $r -le 15; $r++) {
    $wsSheet1.Cells.Item($r, 1).Value = "AnonymousUserTestCases"
}

# Give column A on Sheet1 an explicit width (best-fit sized for the label above)
$wsSheet1.Columns.Item(1).ColumnWidth = 19.83

# Row 12 on Sheet1: fill in Airline / FlightNumber (columns BB / BC)
$wsSheet1.Cells.Item(12, 54).Value = "Delta"
$wsSheet1.Cells.Item(12, 55).Value = 2323

# --- Selections: Anonymous -> B12, Sheet1 -> A16 (select Sheet1 last so it
#     ends up as the active/selected tab, matching the new activeTab) ---
$wsAnonymous.Range("B12").Select() | Out-Null
$wsSheet1.Range("A16").Select() | Out-Null
